$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 15022.4
$ws.Range("I13").Value = 222
$ws.Range("J13").Value = 18722.5
$ws.Range("K13").Value = 222
$ws.Range("L13").Value = 18722.5
$ws.Range("M13").Value = -53
$ws.Range("N13").Value = -19060.5

$ws.Range("H15").Value = 1452.1628
$ws.Range("I15").Value = 1452.1628
$ws.Range("K15").Value = 4356.4884
$ws.Range("M15").Value = -4187.4884

$ws.Range("H18").Value = 969.2143
$ws.Range("I18").Value = 880.75
$ws.Range("K18").Value = 880.75
$ws.Range("M18").Value = -596.75

$ws.Range("H74").Value = 3118.9048
$ws.Range("I74").Value = 3029.1177
$ws.Range("J74").Value = 3500.5
$ws.Range("K74").Value = 3029.1177
$ws.Range("L74").Value = 3500.5
$ws.Range("M74").Value = -2093.1177
$ws.Range("N74").Value = -5372.5

$ws.Range("H77").Value = 3118.9048
$ws.Range("I77").Value = 3029.1177
$ws.Range("J77").Value = 3500.5
$ws.Range("K77").Value = 15145.5885
$ws.Range("L77").Value = 17502.5
$ws.Range("M77").Value = -10465.5885
$ws.Range("N77").Value = -26862.5

$ws.Range("H80").Value = 2050.0952
$ws.Range("I80").Value = 2684.4
$ws.Range("J80").Value = 1851.875
$ws.Range("K80").Value = 8053.200000000001
$ws.Range("L80").Value = 5555.625
$ws.Range("M80").Value = -7055.200000000001
$ws.Range("N80").Value = -7551.625

$ws.Range("H83").Value = 2050.0952
$ws.Range("I83").Value = 2684.4
$ws.Range("J83").Value = 1851.875
$ws.Range("K83").Value = 24159.6
$ws.Range("L83").Value = 16666.875
$ws.Range("M83").Value = -19167.6
$ws.Range("N83").Value = -26650.875

$ws.Range("H88").Value = 24080.893
$ws.Range("I88").Value = 57786
$ws.Range("J88").Value = 5355.8335
$ws.Range("K88").Value = 57786
$ws.Range("L88").Value = 5355.8335
$ws.Range("M88").Value = -57380
$ws.Range("N88").Value = -6167.8335

$ws.Range("H91").Value = 24080.893
$ws.Range("I91").Value = 57786
$ws.Range("J91").Value = 5355.8335
$ws.Range("K91").Value = 57786
$ws.Range("L91").Value = 5355.8335
$ws.Range("M91").Value = -56382
$ws.Range("N91").Value = -8163.8335

$ws.Range("H92").Value = 538.5185
$ws.Range("I92").Value = 590.58826
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 590.58826
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 657.41174
$ws.Range("N92").Value = -2946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H32").Value = 6584.75
$ws.Range("I32").Value = 4487.3037
$ws.Range("J32").Value = 19330.77
$ws.Range("K32").Value = 4487.3037
$ws.Range("L32").Value = 19330.77
$ws.Range("M32").Value = -4200.3037
$ws.Range("N32").Value = -19904.77

$ws.Range("H43").Value = 8315.4
$ws.Range("J43").Value = 8315.4
$ws.Range("L43").Value = 8315.4
$ws.Range("N43").Value = -8941.4

$ws.Range("H61").Value = 1258.0883
$ws.Range("I61").Value = 1102.5862
$ws.Range("J61").Value = 2160
$ws.Range("K61").Value = 1102.5862
$ws.Range("L61").Value = 2160
$ws.Range("M61").Value = -890.5862
$ws.Range("N61").Value = -2584

$ws.Range("H132").Value = 2041.2
$ws.Range("I132").Value = 1301.85
$ws.Range("J132").Value = 4998.6
$ws.Range("K132").Value = 3905.55
$ws.Range("L132").Value = 14995.8
$ws.Range("M132").Value = -1375.55
$ws.Range("N132").Value = -20055.8

$ws.Range("H136").Value = 1258.0883
$ws.Range("I136").Value = 1102.5862
$ws.Range("J136").Value = 2160
$ws.Range("K136").Value = 3307.7586
$ws.Range("L136").Value = 6480
$ws.Range("M136").Value = -757.7586000000001
$ws.Range("N136").Value = -11580

$ws.Range("H141").Value = 39214.5
$ws.Range("J141").Value = 39214.5
$ws.Range("L141").Value = 39214.5
$ws.Range("N141").Value = -49574.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 25000
$ws.Range("J16").Value = 25000
$ws.Range("L16").Value = 25000
$ws.Range("N16").Value = -25340

$ws.Range("H64").Value = 364.72726
$ws.Range("I64").Value = 311.8
$ws.Range("J64").Value = 380.29413
$ws.Range("K64").Value = 311.8
$ws.Range("L64").Value = 380.29413
$ws.Range("M64").Value = -86.80000000000001
$ws.Range("N64").Value = -830.29413

$ws.Range("H67").Value = 364.72726
$ws.Range("I67").Value = 311.8
$ws.Range("J67").Value = 380.29413
$ws.Range("K67").Value = 311.8
$ws.Range("L67").Value = 380.29413
$ws.Range("M67").Value = 468.2
$ws.Range("N67").Value = -1940.29413

$ws.Range("H94").Value = 1061.5116
$ws.Range("I94").Value = 788.0526
$ws.Range("K94").Value = 788.0526
$ws.Range("M94").Value = -337.0526

$ws.Range("H134").Value = 836781.8
$ws.Range("I134").Value = 912204.9
$ws.Range("J134").Value = 7128.5
$ws.Range("K134").Value = 2736614.7
$ws.Range("L134").Value = 21385.5
$ws.Range("M134").Value = -2734079.7
$ws.Range("N134").Value = -26455.5

$ws.Range("H140").Value = 19663.076
$ws.Range("J140").Value = 19663.076
$ws.Range("L140").Value = 19663.076
$ws.Range("N140").Value = -30023.076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 80.55556
$ws.Range("I7").Value = 34.444443
$ws.Range("J7").Value = 126.666664
$ws.Range("K7").Value = 34.444443
$ws.Range("L7").Value = 126.666664
$ws.Range("M7").Value = 78.55555699999999
$ws.Range("N7").Value = -352.666664

$ws.Range("H25").Value = 13175
$ws.Range("I25").Value = 6350
$ws.Range("K25").Value = 6350
$ws.Range("M25").Value = -6176

$ws.Range("H74").Value = 19600
$ws.Range("J74").Value = 22000
$ws.Range("L74").Value = 22000
$ws.Range("N74").Value = -23748

$ws.Range("H77").Value = 19600
$ws.Range("J77").Value = 22000
$ws.Range("L77").Value = 66000
$ws.Range("N77").Value = -74736

$ws.Range("H132").Value = 1544987.8
$ws.Range("I132").Value = 1459.5385
$ws.Range("J132").Value = 3369157.5
$ws.Range("K132").Value = 4378.6155
$ws.Range("L132").Value = 10107472.5
$ws.Range("M132").Value = -1848.6155
$ws.Range("N132").Value = -10112532.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2500
$ws.Range("J20").Value = 4200
$ws.Range("L20").Value = 12600
$ws.Range("N20").Value = -13054

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 15632.667
$ws.Range("I6").Value = 4900
$ws.Range("J6").Value = 20999
$ws.Range("K6").Value = 4900
$ws.Range("L6").Value = 20999
$ws.Range("M6").Value = -4787
$ws.Range("N6").Value = -21225

$ws.Range("H16").Value = 15632.667
$ws.Range("I16").Value = 4900
$ws.Range("J16").Value = 20999
$ws.Range("K16").Value = 4900
$ws.Range("L16").Value = 20999
$ws.Range("M16").Value = -4650
$ws.Range("N16").Value = -21499

$ws.Range("H22").Value = 5263
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 5263
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 5263
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6321

$ws.Range("H113").Value = 1607
$ws.Range("I113").Value = 1358.7142
$ws.Range("J113").Value = 1896.6666
$ws.Range("K113").Value = 1358.7142
$ws.Range("L113").Value = 1896.6666
$ws.Range("M113").Value = 811.2858000000001
$ws.Range("N113").Value = -6236.6666

$ws.Range("H132").Value = 2503073.8
$ws.Range("I132").Value = 2892.4546
$ws.Range("J132").Value = 14289643
$ws.Range("K132").Value = 8677.363799999999
$ws.Range("L132").Value = 42868929
$ws.Range("M132").Value = -6147.363799999999
$ws.Range("N132").Value = -42873989

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 14506.75
$ws.Range("I13").Value = 8006
$ws.Range("J13").Value = 16673.666
$ws.Range("K13").Value = 8006
$ws.Range("L13").Value = 16673.666
$ws.Range("M13").Value = -7866
$ws.Range("N13").Value = -16953.666

$ws.Range("H22").Value = 23230.818
$ws.Range("I22").Value = 55908.11
$ws.Range("J22").Value = 608.0769
$ws.Range("K22").Value = 55908.11
$ws.Range("L22").Value = 608.0769
$ws.Range("M22").Value = -55613.11
$ws.Range("N22").Value = -1198.0769

$ws.Range("H23").Value = 838418.7
$ws.Range("I23").Value = 2502003
$ws.Range("J23").Value = 6626.5
$ws.Range("K23").Value = 2502003
$ws.Range("L23").Value = 6626.5
$ws.Range("M23").Value = -2501773
$ws.Range("N23").Value = -7086.5

$ws.Range("H27").Value = 23230.818
$ws.Range("I27").Value = 55908.11
$ws.Range("J27").Value = 608.0769
$ws.Range("K27").Value = 55908.11
$ws.Range("L27").Value = 608.0769
$ws.Range("M27").Value = -55801.11
$ws.Range("N27").Value = -822.0769

$ws.Range("H46").Value = 808.6786
$ws.Range("J46").Value = 749.5789
$ws.Range("L46").Value = 749.5789
$ws.Range("N46").Value = -1125.5789

$ws.Range("H55").Value = 867.58826
$ws.Range("I55").Value = 417.5
$ws.Range("J55").Value = 1006.0769
$ws.Range("K55").Value = 417.5
$ws.Range("L55").Value = 1006.0769
$ws.Range("M55").Value = -244.5
$ws.Range("N55").Value = -1352.0769

$ws.Range("H138").Value = 24452.666
$ws.Range("J138").Value = 24452.666
$ws.Range("L138").Value = 24452.666
$ws.Range("N138").Value = -34732.666

$ws.Range("H140").Value = 41816.355
$ws.Range("J140").Value = 41816.355
$ws.Range("L140").Value = 41816.355
$ws.Range("N140").Value = -52176.355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1413.2307
$ws.Range("I100").Value = 1252.1111
$ws.Range("J100").Value = 1775.75
$ws.Range("K100").Value = 2504.2222
$ws.Range("L100").Value = 3551.5
$ws.Range("M100").Value = -1963.2222
$ws.Range("N100").Value = -4633.5

